$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1 becomes a proper header row ---
# (it previously held a stray duplicate of row 2's data instead of field names)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Carry the existing header styling (bold font + border) from B1:G1 onto the newly added H1:N1 cells
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 gains the remaining record columns ---
# (B2/C2/D2/E2/F2/G2 already hold the correct values: name/capacity/owner/register_date/register_reason/acquire_value)
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("K2").Value = "蔣乃辛"
$ws.Range("L2").Value = 1722
$ws.Range("M2").Value = "tmp7091"
$ws.Range("N2").Value = 31

# "date" column: stage the literal text via a scratch formula cell and paste values-only,
# so Excel stores the plain string "2012-04-20" instead of reinterpreting it as a date serial.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""2012-04-20"""
$scratch.Copy()
$ws.Range("J2").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$scratch.ClearContents()

# Carry the existing data-row styling from B2:G2 onto the newly added H2:N2 cells
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
